$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of score data appended below the existing table.
# A5 must hold the literal text "2026-02-10" (not an auto-converted date
# serial) to match the existing Date column, which is stored as plain text.
# Entering it as a formula that evaluates to the string, then pasting the
# result back as a value, keeps Excel's "smart" date-detection from firing
# while leaving the cell's style untouched (same as the existing rows).
$ws.Cells.Item(5, 1).Formula = '="2026-02-10"'
$ws.Cells.Item(5, 1).Copy() | Out-Null
$ws.Cells.Item(5, 1).PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = 0

$ws.Cells.Item(5, 2).Value = 865
